# Parameters_DNAs.xlsx edit:
# - Row 8 (DNAscent_T5): update script body, cpus-per-task stays 32, threads 64->10,
#   "no specified" -> "No specified", "Yes" -> "No", times "11H" -> "5H" (both columns)
# - Row 9 (DNAscent_T6): update script body, cpus-per-task 32->10, threads 64->10,
#   RAM 256 -> "128G", times "10H"/"20 H" -> "2 H"/"2H "
# - Row 10 (DNAscent_T7): update script body (now reuses the new T6 script),
#   threads 32->30, "no specified" -> "No specified", comment "" -> "Fastest"
# - Remove old rows 11 (DNAscent_T8) and 12 (DNAscent_T9) entirely
# - Row heights: row 8 153->136, row 9 170->153 (row 10 stays 153)
# - Selection moves to K10, viewport scrolled so row 7 is visible at top

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: DNAscent_T5 ---
$ws.Range("B8").Value = "#! /usr/bin/bash`n#SBATCH --job-name=T5`n#SBATCH --ntasks=1`n#SBATCH --cpus-per-task=32`n#SBATCH --time=72:00:0`n#SBATCH --partition=ncpu`n#SBATCH --output=std/dnascent_T1.o`n#SBATCH --error=std/dnascent_T1.e"
$ws.Range("D8").Value = 10
$ws.Range("E8").Value = "No specified"
$ws.Range("F8").Value = "No"
$ws.Range("G8").Value = "5H"
$ws.Range("H8").Value = "5H"
$ws.Rows.Item(8).RowHeight = 136

# --- Row 9: DNAscent_T6 ---
$ws.Range("B9").Value = "#! /usr/bin/bash`n#SBATCH --job-name=T6`n#SBATCH --ntasks=1`n#SBATCH --cpus-per-task=12`n#SBATCH --gres=gpu:1`n#SBATCH --time=16:00:0`n#SBATCH --partition=gpu`n#SBATCH --output=std/dnascent_T6.o`n#SBATCH --error=std/dnascent_T6.e"
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = "128G"
$ws.Range("G9").Value = "2 H"
$ws.Range("H9").Value = "2H "
$ws.Rows.Item(9).RowHeight = 153

# --- Row 10: DNAscent_T7 ---
$ws.Range("B10").Value = "#! /usr/bin/bash`n#SBATCH --job-name=T6`n#SBATCH --ntasks=1`n#SBATCH --cpus-per-task=12`n#SBATCH --gres=gpu:1`n#SBATCH --time=16:00:0`n#SBATCH --partition=gpu`n#SBATCH --output=std/dnascent_T6.o`n#SBATCH --error=std/dnascent_T6.e"
$ws.Range("D10").Value = 30
$ws.Range("E10").Value = "No specified"
$ws.Range("I10").Value = "Fastest"

# --- Remove old rows 11 (DNAscent_T8) and 12 (DNAscent_T9) ---
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(11).Delete()

# --- View state: selection + scroll position ---
$ws.Range("K10").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
